# Insert a new price-report row at row 164 (weekly Betarraga data point),
# pushing the existing rows 164..252 down to 165..253.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(164).Insert()

$ws.Cells.Item(164, 1).Value = 7
$ws.Cells.Item(164, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(164, 3).Value = "Ñuble"
$ws.Cells.Item(164, 4).Value = 44455
$ws.Cells.Item(164, 5).Value = 16
$ws.Cells.Item(164, 6).Value = 100114014
$ws.Cells.Item(164, 7).Value = "Betarraga"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 300
$ws.Cells.Item(164, 11).Value = 750
$ws.Cells.Item(164, 12).Value = 800
$ws.Cells.Item(164, 13).Value = 775
$ws.Cells.Item(164, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(164, 15).Value = "Región del Maule"
$ws.Cells.Item(164, 16).Value = 155
$ws.Cells.Item(164, 17).Value = 5
$ws.Cells.Item(164, 18).Value = "Hortaliza"
